# Updated test data for 5,24,40V,BatteryStandby and AC Calculations test cases
#
# The only functional content change in the workbook is the Device name in
# cell G8 of the "Add Panels and Devices" sheet: "PRN800" -> "PX-PR".
# (All other differences in the target OOXML are shared-string re-indexing
# that falls out automatically once the string table no longer needs the
# now-unused "PRN800" entry, plus Excel-version/session metadata bumps that
# are not meaningful COM operations.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G8").Value = "PX-PR"
